# Append the final data-output row (row 12) of sequential execution times
# for the K and L columns (the two columns whose measurements completed
# after everything else had already been written out).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(12, 11).Value = 0.6609286   # K12
$ws.Cells.Item(12, 12).Value = 0.5932979   # L12
